# "change name shipping option" - rename the shipping-option labels used in the
# shippingOption column on the FPA sheets. The old labels ("Standard shipping",
# "Pickup at factory", "Special packaging / via freight forwarding") are
# replaced by new ones. Writing the "Freight delivery ..." value first (then
# "Package delivery ...", then "Pick-up at factory ...") keeps the shared
# string insertion order identical across all sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FPA011")
$ws1.Range("P4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws1.Range("P2").Value = "Package delivery (extra costs)"
$ws1.Range("P3").Value = "Pick-up at factory (no costs)"
$ws1.Range("P5").Value = "Package delivery (extra costs)"

$ws2 = $wb.Worksheets.Item("FPA012-013-015-017")
$ws2.Range("P4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws2.Range("P2").Value = "Package delivery (extra costs)"
$ws2.Range("P3").Value = "Pick-up at factory (no costs)"
$ws2.Range("P5").Value = "Package delivery (extra costs)"

$ws3 = $wb.Worksheets.Item("FPA014-016-020")
$ws3.Range("P2").Value = "Package delivery (extra costs)"
$ws3.Range("P3").Value = "Pick-up at factory (no costs)"

$ws4 = $wb.Worksheets.Item("FPA018-019")
$ws4.Range("Q4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws4.Range("Q2").Value = "Package delivery (extra costs)"
$ws4.Range("Q3").Value = "Pick-up at factory (no costs)"
$ws4.Range("Q5").Value = "Package delivery (extra costs)"

# Active-sheet / selection bookkeeping: FPA011 becomes the selected tab
# (instead of FPA012-013-015-017), with M19 selected there.
$ws1.Activate() | Out-Null
$ws1.Range("M19").Select() | Out-Null
